# ENWIAM test-case worksheet update
# - Row 26 (ENWIAM41) description/jira text updated, row height 45 -> 30
# - Row 27 (ENWIAM42) left as-is
# - New row 28 (ENWIAM43) inserted
# - Old rows 28/29 (ENWIAM50/ENWIAM51) shift down to 29/30, content unchanged
# - sheet selection moved to D28

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update row 26 content in place (TCID stays ENWIAM41) ---
$ws.Range("B26").Value = " OPQA-1859"
$ws.Range("C26").Value = "Verify that user who signed in to Neon through social shall be able to navigate to ENW after providing steam password in the Linking Modal"

# --- Step 2: insert a brand-new row at position 28 ---
# This pushes the current row 28 (ENWIAM50) -> 29 and row 29 (ENWIAM51) -> 30,
# automatically preserving their contents/styles.
$ws.Rows.Item(28).Insert()

# Copy the formatting from row 27 (same visual style family as row 26/27)
# onto the freshly inserted (blank) row 28.
$ws.Range("A27:E27").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: populate the new row 28 (ENWIAM43) ---
$ws.Range("A28").Value = "ENWIAM43"
$ws.Range("B28").Value = "OPQA-1686"
$ws.Range("C28").Value = 'Verify that the User is able to see message  "Your account registration has not yet been confirmed. Please click on the verification link you were sent by email from <no-reply-email-address>, or have a link resent.'
$ws.Range("D28").Value = "Y"

# --- Step 4: row heights ---
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 45

# --- Step 5: update selection / view ---
[void]$ws.Activate()
[void]$ws.Range("D28").Select()

Write-Host "ENWIAM sheet updated"
